$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently sits right
#    after the title heading ("Play Always Hot Deluxe Slot Game for
#    Free - Review").
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------
# 2. Insert a new bold paragraph ("Play Always Hot Deluxe Slot Game for
#    Free - Review") right before the final paragraph of the document
#    (the one that currently holds the "Create a feature image..."
#    image-prompt text).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$startOfLast = $lastPara.Range.Start
$insertionPoint = $d.Range($startOfLast, $startOfLast)

$headingText = "Play Always Hot Deluxe Slot Game for Free - Review"
$insertionPoint.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>" + $headingText + "</w:t></w:r></w:p>")

# Split the freshly inserted run from the original final paragraph so
# that they become two distinct paragraphs.
$boundary = $startOfLast + $headingText.Length
$splitPoint = $d.Range($boundary, $boundary)
$splitPoint.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 3. Replace the body text of what is now the last paragraph (formerly
#    the "Create a feature image..." prompt) with the meta-description
#    copy, keeping its existing italic formatting intact.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($count)

$oldText = "Create a feature image that captures the essence of Always Hot Deluxe. The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a slot machine with flames and hot red fruits bursting out of it. The background should be filled with flames and the words " + [char]34 + "Always Hot Deluxe" + [char]34 + " should be written in bold and fiery letters."
$newText = "Read our unbiased review of Always Hot Deluxe, a classic slot machine with a non-progressive jackpot of up to 60,000 coins. Play for free and learn more!"

$finalPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
